$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Example")

# -----------------------------------------------------------------
# 1) Stage formatting templates from the existing rows into scratch
#    rows far below the table, so we can freely overwrite the
#    original rows afterwards without losing the formats we need.
# -----------------------------------------------------------------
$ws.Range("A2:E2").Copy($ws.Range("A50:E50"))   # "plain" style  (white, no fill)
$ws.Range("A3:E3").Copy($ws.Range("A51:E51"))   # "plain" style  (white, no fill)
$ws.Range("A6:E6").Copy($ws.Range("A52:E52"))   # "highlighted, default font" style
$ws.Range("A8:E8").Copy($ws.Range("A53:E53"))   # "highlighted, bigger font (E col)" style
$ws.Range("A9:E9").Copy($ws.Range("A54:E54"))   # "highlighted, bigger font (E col)" style

# -----------------------------------------------------------------
# 2) Remove the old merged note cells (F8:G10 / F11:G13) and the
#    F/G columns entirely (no longer used in the new table).
# -----------------------------------------------------------------
$ws.Range("F8:G10").UnMerge()
$ws.Range("F11:G13").UnMerge()
$ws.Columns("F:G").Delete()

# -----------------------------------------------------------------
# 3) Drop the old trailing rows (11-13) that are no longer part of
#    the new, shorter table (A1:E10).
# -----------------------------------------------------------------
$ws.Rows("11:13").Delete()

# -----------------------------------------------------------------
# 4) Rebuild rows 2-10 from the staged templates (this brings the
#    correct cell formatting along with the copy), then overwrite
#    the cell values with the new example data.
# -----------------------------------------------------------------
$ws.Range("A53:E53").Copy($ws.Range("A2:E2"))
$ws.Range("A54:E54").Copy($ws.Range("A3:E3"))
$ws.Range("A52:E52").Copy($ws.Range("A4:E4"))
$ws.Range("A50:E50").Copy($ws.Range("A5:E5"))
$ws.Range("A51:E51").Copy($ws.Range("A6:E6"))
$ws.Range("A52:E52").Copy($ws.Range("A7:E7"))
$ws.Range("A52:E52").Copy($ws.Range("A8:E8"))
$ws.Range("A52:E52").Copy($ws.Range("A9:E9"))
$ws.Range("A52:E52").Copy($ws.Range("A10:E10"))

# Row 4 is the old "highlighted" P002 row with the highlight
# switched off (explicit "No Fill") - matches the workbook's own
# pattern of re-using a highlighted row's font but clearing the fill.
$ws.Range("A4:E4").Interior.ColorIndex = -4142   # xlColorIndexNone / "No Fill"

# -----------------------------------------------------------------
# 5) Clear the scratch rows used for staging.
# -----------------------------------------------------------------
$ws.Rows("50:54").Delete()

# -----------------------------------------------------------------
# 6) Write the new example values.
# -----------------------------------------------------------------
# Row 2: P001, 44318, 1, (blank), DIP-PER-TET
$ws.Range("A2").Value = "P001"
$ws.Range("B2").Value = 44318
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "DIP-PER-TET"

# Row 3: P001, 44404, 2, (blank), DIP-PER-TET
$ws.Range("A3").Value = "P001"
$ws.Range("B3").Value = 44404
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "DIP-PER-TET"

# Row 4: P002, 44428, 1, pfizer, COV
$ws.Range("A4").Value = "P002"
$ws.Range("B4").Value = 44428
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = "pfizer"
$ws.Range("E4").Value = "COV"

# Row 5: P002, 44407, 2, pfizer, COV
$ws.Range("A5").Value = "P002"
$ws.Range("B5").Value = 44407
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = "pfizer"
$ws.Range("E5").Value = "COV"

# Row 6: P002, 44772, 1, (blank), HPV
$ws.Range("A6").Value = "P002"
$ws.Range("B6").Value = 44772
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = "HPV"

# Row 7: P003, 44197, 1, astrazeneca, COV
$ws.Range("A7").Value = "P003"
$ws.Range("B7").Value = 44197
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = "astrazeneca"
$ws.Range("E7").Value = "COV"

# Row 8: P003, 44221, 2, astrazeneca, COV
$ws.Range("A8").Value = "P003"
$ws.Range("B8").Value = 44221
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = "astrazeneca"
$ws.Range("E8").Value = "COV"

# Row 9: P003, 44287, 3, pfizer, COV
$ws.Range("A9").Value = "P003"
$ws.Range("B9").Value = 44287
$ws.Range("C9").Value = 3
$ws.Range("D9").Value = "pfizer"
$ws.Range("E9").Value = "COV"

# Row 10: P003, 44481, 4, moderna, COV
$ws.Range("A10").Value = "P003"
$ws.Range("B10").Value = 44481
$ws.Range("C10").Value = 4
$ws.Range("D10").Value = "moderna"
$ws.Range("E10").Value = "COV"

# -----------------------------------------------------------------
# 7) Cosmetic touch-ups to line up with the final layout.
# -----------------------------------------------------------------
$ws.Columns("E").ColumnWidth = 12.11
$ws.Range("C10").Select()

Write-Output "done"
